$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update the computed values (E2/F2/F5) ---
$ws.Range("E2").Value = 0.2026
$ws.Range("F2").Value = 0.041
$ws.Range("F5").Value = 0.084

# --- Reapply formatting to force Excel's style-table dedup/reorder ---
$ws.Range("A2:A5").Interior.ThemeColor = 5
$ws.Range("A7:A9").Interior.ThemeColor = 5
$ws.Range("B6:F6").Borders.Item(8).LineStyle = 1

# --- Update view state: scroll so column C is the left-most visible column,
#     and select G15 as the active cell ---
$ws.Range("G15").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
